$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: re-label + reorder columns, drop the old H1/I1 "fout op functiegenerator freq" /
#     "fout op stroom" labels that no longer apply ---
$ws.Range("A1").Value = "current"
$ws.Range("B1").Value = "current_err"
$ws.Range("C1").Value = "meting"
$ws.Range("H1").Value = $null
$ws.Range("I1").Value = $null

# --- Data rows: keep A (current) / B (current_err formula) / C (meting index) in sync.
#     Row 6's old text annotation ("fout op stroom") becomes a plain sequential index, and an
#     extra duplicate measurement row is inserted right after it (A repeats 0.249, meting 5->6).
#     Likewise row 9's old text annotation ("9, 10") becomes a plain sequential index, with another
#     duplicate measurement row inserted right after it (A repeats 0.401, meting 9->10).
#     Every row below that shifts down by two and all subsequent measurement indices renumber
#     sequentially through 17. ---

$aValues = @(0.049, 0.1, 0.15, 0.199, 0.249, 0.249, 0.301, 0.351, 0.401, 0.401, 0.451, 0.499, 0.549, 0.6, 0.649, 0.7, 0.749)
$cValues = @(1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17)

for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $aValues[$i]
    $ws.Range("B$row").Formula = "=IF(A$row<0.4, A$row * 0.016 + 4 * 0.00001, A$row * 0.02 + 10 * 0.001)"
    $ws.Range("C$row").Value = $cValues[$i]
}

$ws.Range("B2:B18").Select()
